$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
# Copy H1's formatting (bold font, thin border, centered/top alignment)
# to the new header cells before setting their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-10
$i0 = @(6, 1, 5, 4, 4, 4, 7, 7, 7)
$if = @(8, 2, 6, 4, 5, 4, 7, 8, 7)

for ($r = 0; $r -lt 9; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
